$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": append row 63 ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Cells.Item(63,1).NumberFormat = $ws1.Cells.Item(62,1).NumberFormat
$ws1.Cells.Item(63,1).Value2 = 45676.99999999999
$ws1.Cells.Item(63,2).Value2 = 48

# --- Sheet "Monthly Trend": append row 22 ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Cells.Item(22,1).NumberFormat = $ws2.Cells.Item(21,1).NumberFormat
$ws2.Cells.Item(22,1).Value2 = 45688.99999999999
$ws2.Cells.Item(22,2).Value2 = 48

# --- Sheet "PO Forecast": new forecast model ---
$ws3 = $wb.Worksheets.Item("PO Forecast")

$ws3.Cells.Item(2,2).Value2 = 620
$ws3.Cells.Item(3,2).Value2 = 20
$ws3.Cells.Item(4,2).Value2 = 0
$ws3.Cells.Item(5,2).Value2 = 0
$ws3.Cells.Item(6,2).Value2 = 0
$ws3.Cells.Item(7,2).Value2 = 0
$ws3.Cells.Item(8,2).Value2 = 14
$ws3.Cells.Item(9,2).Value2 = 0
$ws3.Cells.Item(10,2).Value2 = 0
$ws3.Cells.Item(11,2).Value2 = 30
$ws3.Cells.Item(12,2).Value2 = 196
$ws3.Cells.Item(13,2).Value2 = 376
$ws3.Cells.Item(14,2).Value2 = 456
$ws3.Cells.Item(15,2).Value2 = 395
$ws3.Cells.Item(16,2).Value2 = 239
$ws3.Cells.Item(17,2).Value2 = 199
$ws3.Cells.Item(18,2).Value2 = 108
$ws3.Cells.Item(19,2).Value2 = 195
$ws3.Cells.Item(20,2).Value2 = 123
$ws3.Cells.Item(21,2).Value2 = 0
$ws3.Cells.Item(22,2).Value2 = 32
$ws3.Cells.Item(23,2).Value2 = 184
$ws3.Cells.Item(24,2).Value2 = 224
$ws3.Cells.Item(25,2).Value2 = 148
$ws3.Cells.Item(26,2).Value2 = 122
$ws3.Cells.Item(27,2).Value2 = 266
$ws3.Cells.Item(28,2).Value2 = 468
$ws3.Cells.Item(29,2).Value2 = 354
$ws3.Cells.Item(30,2).Value2 = 334
$ws3.Cells.Item(31,2).Value2 = 802
$ws3.Cells.Item(32,2).Value2 = 1058
$ws3.Cells.Item(33,2).Value2 = 0
$ws3.Cells.Item(34,2).Value2 = 89
$ws3.Cells.Item(35,2).Value2 = 17
$ws3.Cells.Item(36,2).Value2 = 40
$ws3.Cells.Item(37,2).Value2 = 163
$ws3.Cells.Item(38,2).Value2 = 268
$ws3.Cells.Item(39,2).Value2 = 266
$ws3.Cells.Item(40,2).Value2 = 191
$ws3.Cells.Item(41,2).Value2 = 137
$ws3.Cells.Item(42,2).Value2 = 149
$ws3.Cells.Item(43,2).Value2 = 203
$ws3.Cells.Item(44,2).Value2 = 270
$ws3.Cells.Item(45,2).Value2 = 259
$ws3.Cells.Item(46,2).Value2 = 227
$ws3.Cells.Item(47,2).Value2 = 210
$ws3.Cells.Item(48,2).Value2 = 265
$ws3.Cells.Item(49,2).Value2 = 418
$ws3.Cells.Item(50,2).Value2 = 605
$ws3.Cells.Item(51,2).Value2 = 709
$ws3.Cells.Item(52,2).Value2 = 670
$ws3.Cells.Item(53,2).Value2 = 476
$ws3.Cells.Item(54,2).Value2 = 488
$ws3.Cells.Item(55,2).Value2 = 377
$ws3.Cells.Item(56,2).Value2 = 338
$ws3.Cells.Item(57,2).Value2 = 450
$ws3.Cells.Item(58,2).Value2 = 401
$ws3.Cells.Item(59,2).Value2 = 262
$ws3.Cells.Item(60,2).Value2 = 261
$ws3.Cells.Item(61,2).Value2 = 418
$ws3.Cells.Item(62,2).Value2 = 694
$ws3.Cells.Item(63,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(63,1).Value2 = 45676.99999999999
$ws3.Cells.Item(63,2).Value2 = 185
$ws3.Cells.Item(64,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(64,1).Value2 = 45683.99999999999
$ws3.Cells.Item(64,2).Value2 = 223
$ws3.Cells.Item(65,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(65,1).Value2 = 45690.99999999999
$ws3.Cells.Item(65,2).Value2 = 354
$ws3.Cells.Item(66,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(66,1).Value2 = 45697.99999999999
$ws3.Cells.Item(66,2).Value2 = 357
$ws3.Cells.Item(67,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(67,1).Value2 = 45704.99999999999
$ws3.Cells.Item(67,2).Value2 = 282
$ws3.Cells.Item(68,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(68,1).Value2 = 45711.99999999999
$ws3.Cells.Item(68,2).Value2 = 282
$ws3.Cells.Item(69,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(69,1).Value2 = 45718.99999999999
$ws3.Cells.Item(69,2).Value2 = 395
$ws3.Cells.Item(70,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(70,1).Value2 = 45725.99999999999
$ws3.Cells.Item(70,2).Value2 = 512
$ws3.Cells.Item(71,1).NumberFormat = $ws3.Cells.Item(62,1).NumberFormat
$ws3.Cells.Item(71,1).Value2 = 45732.99999999999
$ws3.Cells.Item(71,2).Value2 = 530
